# Update the cryptocurrency price/volume table (rows 2-51) to match
# the latest scraped snapshot. Column D ("Price") values are forced to
# stay plain text (matching the original inline-string cells, which use
# "." as a thousands separator and carry significant trailing zeros -
# Excel would otherwise auto-convert them to numbers and mangle them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.421.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4471"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3761"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8932"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.826.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.754"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.410"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07107"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008817"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.432.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.271"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.052.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.983"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.381"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.370"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08836"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7833"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.205"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.530"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.111"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01991"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05331"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.387"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1732"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.865"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.296"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +18.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.797"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("E46").Value = "  +8.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.703"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06375"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
